# Weekly update: insert a new week's worth of "Repollo" price rows
# (Terminal La Palmera de La Serena) at rows 973-978, shifting the
# existing data down by 6 rows (so the sheet grows from 1073 to 1079
# data rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 blank rows at the top of the block (this pushes the old
# rows 973:1073 down to 979:1079).
$ws.Rows("973:978").Insert()

# Fixed/common column values shared by every "Repollo" row in this sheet.
$mercadoId   = 8
$mercado     = "Terminal La Palmera de La Serena"
$region      = "Coquimbo"
$codreg      = 4
$categoriaId = 100112006
$categoria   = "Repollo"
$unidad      = "$/unidad"
$origen      = "Provincia del Elquí"
$kgUnidades  = 1
$clasif      = "Hortaliza"

# New week's data: Fecha (serial), Variedad, Calidad, Volumen, Precio
# minimo, Precio maximo, Precio promedio ponderado (= Precio $/Kg).
$newRows = @(
    @{ Row = 973; Fecha = 44449; Variedad = "Copenhague";    Calidad = "Primera"; Volumen = 2800; Min = 850; Max = 900; Prom = 875 },
    @{ Row = 974; Fecha = 44449; Variedad = "Copenhague";    Calidad = "Segunda"; Volumen = 1400; Min = 750; Max = 800; Prom = 775 },
    @{ Row = 975; Fecha = 44449; Variedad = "Crespo record"; Calidad = "Primera"; Volumen = 2600; Min = 700; Max = 800; Prom = 750 },
    @{ Row = 976; Fecha = 44449; Variedad = "Crespo record"; Calidad = "Segunda"; Volumen = 1400; Min = 500; Max = 600; Prom = 550 },
    @{ Row = 977; Fecha = 44449; Variedad = "Morada(o)";     Calidad = "Primera"; Volumen = 2000; Min = 700; Max = 800; Prom = 750 },
    @{ Row = 978; Fecha = 44449; Variedad = "Morada(o)";     Calidad = "Segunda"; Volumen = 1000; Min = 500; Max = 600; Prom = 550 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = $mercadoId
    $ws.Cells.Item($row, 2).Value  = $mercado
    $ws.Cells.Item($row, 3).Value  = $region
    $ws.Cells.Item($row, 4).Value  = $r.Fecha
    $ws.Cells.Item($row, 5).Value  = $codreg
    $ws.Cells.Item($row, 6).Value  = $categoriaId
    $ws.Cells.Item($row, 7).Value  = $categoria
    $ws.Cells.Item($row, 8).Value  = $r.Variedad
    $ws.Cells.Item($row, 9).Value  = $r.Calidad
    $ws.Cells.Item($row, 10).Value = $r.Volumen
    $ws.Cells.Item($row, 11).Value = $r.Min
    $ws.Cells.Item($row, 12).Value = $r.Max
    $ws.Cells.Item($row, 13).Value = $r.Prom
    $ws.Cells.Item($row, 14).Value = $unidad
    $ws.Cells.Item($row, 15).Value = $origen
    $ws.Cells.Item($row, 16).Value = $r.Prom
    $ws.Cells.Item($row, 17).Value = $kgUnidades
    $ws.Cells.Item($row, 18).Value = $clasif
}
